# "Settings back to qfin22 presentation - recalibration"
# Recalibrated numeric outputs across the AR / SETAR / GARCH / TARCH / AR_TARCH
# sheets, plus the autocorrelation-vector strings stored alongside them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.9335904046783494
$ws.Range("B3").Value = 0.8054165087168166
$ws.Range("B4").Value = 53093.55241370312
$ws.Range("B5").Value = "[1.0, 0.3584453424327151, 0.37930025599873984, 0.3510771181268938, 0.335221735811094, 0.5356913162204199, 0.23728150865475217, 0.28725459481259835, 0.24998707665642506, 0.24648087882497297, 0.2402098674124782, 0.25970392829889355, 0.2688381466271559, 0.2638332863759762, 0.2663104083545547, 0.25820451953065354, 0.2749460172287829, 0.26189166406091174, 0.24618896377493396, 0.2625179666369773]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -25.7716308955992
$ws.Range("B5").Value = 0.8150208989744211
$ws.Range("B6").Value = 72890.23567939775
$ws.Range("B7").Value = 20.50126654294485
$ws.Range("B8").Value = 0.8239600773612403
$ws.Range("B9").Value = 19471.53473012591
$ws.Range("B10").Value = "[1.0, 0.25213620269477477, 0.26134641640547895, 0.25330501644533654, 0.27716337765534316, 0.2603351933874866, 0.19531207702656486, 0.2522354391607952, 0.44812358600260743, 0.16716026892441188, 0.1813303399939047, 0.1828774983281022, 0.18553762093053155, 0.18446367684638504, 0.19437580594903203, 0.19470646957668575, 0.2000152509090146, 0.18226707284066837, 0.16929983480764305, 0.18403642452132976]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.003833604654296428
$ws.Range("B3").Value = 0.200651088789857
$ws.Range("B4").Value = 0.1723440371158349
$ws.Range("B5").Value = 0.8276560125101694
$ws.Range("B6").Value = "[1.0, -0.06769554961986754, -0.021270611760886023, -0.055551138151386056, -0.10134888106358124, 0.2447878613960221, -0.0345535495498049, 0.04343116294860986, 0.021587830889403577, -0.04339030192318453, -0.055010180283633775, -0.0387778554029656, -0.005578411045380245, 0.024560103136938415, 0.035763976497856016, 0.00519301133682494, -0.011197722935551954, -0.028458375071271136, -0.018520779558156883, 0.01820803589320333]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.09503409356902032
$ws.Range("B3").Value = 0.2414968174470326
$ws.Range("B4").Value = 0.1040083720951872
$ws.Range("B5").Value = 0.8184187314247282
$ws.Range("B6").Value = "[1.0, -0.06444024912362192, -0.020955216777585068, -0.053865392406555766, -0.09923453127764421, 0.2473851303587318, -0.037608816529625916, 0.039751540735715546, 0.022015506742874995, -0.04615675187183926, -0.0567236447929427, -0.0399447574818361, -0.005714901517733225, 0.024593723414622595, 0.035600946577843705, 0.00695609280847492, -0.010776730242895282, -0.025488118335514807, -0.018478019367355165, 0.022381774258670592]"
$ws.Range("B7").Value = 0.1551482674415018

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = -0.00469499547455669
$ws.Range("B3").Value = 0.168186338742623
$ws.Range("B4").Value = 0.1545037286681763
$ws.Range("B5").Value = 0.8356458734655779
$ws.Range("B6").Value = "[1.0, -0.04530717739393527, -0.012199268730498106, -0.046184611898848246, -0.06675224375725772, 0.18165127664149253, -0.033599963012538266, 0.026680636043106314, 0.004060321154986111, -0.0405204665182707, -0.03726272973584698, -0.02938560324023415, -0.019407014705202834, 0.030034254889277797, 0.028678311106904067, 0.0027082489375989537, -0.004082263094791631, -0.019461228851109628, -0.00381435239057834, 0.01775468732643011]"
$ws.Range("B7").Value = 0.01970464452165983
$ws.Range("B9").Value = 0.8094833226602565
